# Generate Report for Handoff
#
# The "dc147130-2548-44bd-a7ed-c8fc363019ce.md" file has finished translation
# and is now ready for handoff. Update its status (and the corresponding
# handoff timestamps) across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: summary row for dc147130-...-c8fc363019ce.md (row 3) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-03-22 12:23:27"

# --- zh-cn detail sheet: same file's row (row 3) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "2016-03-22 12:23:24"

# --- de-de detail sheet: same file's row (row 3) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "2016-03-22 12:23:27"
